# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the newly generated numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 697
$ws1.Range("F3").Value = 15
$ws1.Range("F7").Value = 54
$ws1.Range("F9").Value = 4334
$ws1.Range("F10").Value = 4307
$ws1.Range("F11").Value = 11

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 697
$ws4.Range("F3").Value = 15
$ws4.Range("F7").Value = 54
$ws4.Range("F9").Value = 4335
$ws4.Range("F10").Value = 4307
$ws4.Range("F11").Value = 11
